$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E (particip) and F (taxa_sucesso) for rows 2-7 were stored as
# fractions (e.g. 0.9086...) and are now stored as the equivalent
# percentage number (e.g. 90.86...), while keeping the 0.00% display
# format. Multiply each value by 100 in place.
for ($row = 2; $row -le 7; $row++) {
    foreach ($col in @("E", "F")) {
        $cell = $ws.Range("$col$row")
        $cell.Value2 = $cell.Value2 * 100
    }
}
